$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row update (source hunk #0)
$ws.Range("H135").Value = 5410
$ws.Range("I135").Value = 7049.316
$ws.Range("K135").Value = 63443.844
$ws.Range("M135").Value = -60908.844
# Row update (source hunk #1)
$ws.Range("H138").Value = 348596
$ws.Range("I138").Value = 599318.9399999999
$ws.Range("J138").Value = 3851.9375
$ws.Range("K138").Value = 1797956.82
$ws.Range("L138").Value = 11555.8125
$ws.Range("M138").Value = -1792816.82
$ws.Range("N138").Value = -21835.8125
# Row update (source hunk #2)
$ws.Range("H139").Value = 131072.33
$ws.Range("J139").Value = 131072.33
$ws.Range("L139").Value = 131072.33
$ws.Range("N139").Value = -141352.33

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row update (source hunk #3)
$ws.Range("H61").Value = 18833.834
$ws.Range("I61").Value = 27430
$ws.Range("K61").Value = 27430
$ws.Range("M61").Value = -27218
# Row update (source hunk #4)
$ws.Range("H74").Value = 12278.909
$ws.Range("I74").Value = 14695.875
$ws.Range("J74").Value = 5833.6665
$ws.Range("K74").Value = 14695.875
$ws.Range("L74").Value = 5833.6665
$ws.Range("M74").Value = -13821.875
$ws.Range("N74").Value = -7581.6665
# Row update (source hunk #5)
$ws.Range("H77").Value = 12278.909
$ws.Range("I77").Value = 14695.875
$ws.Range("J77").Value = 5833.6665
$ws.Range("K77").Value = 73479.375
$ws.Range("L77").Value = 29168.3325
$ws.Range("M77").Value = -69111.375
$ws.Range("N77").Value = -37904.3325
# Row update (source hunk #6)
$ws.Range("H96").Value = 34999.5
$ws.Range("J96").Value = 34999.5
$ws.Range("L96").Value = 34999.5
$ws.Range("N96").Value = -40491.5
# Row update (source hunk #7)
$ws.Range("H132").Value = 4085.64
$ws.Range("I132").Value = 3982.3333
$ws.Range("K132").Value = 11946.9999
$ws.Range("M132").Value = -9416.999899999999
# Row update (source hunk #8)
$ws.Range("H136").Value = 18833.834
$ws.Range("I136").Value = 27430
$ws.Range("K136").Value = 82290
$ws.Range("M136").Value = -79740
# Row update (source hunk #9)
$ws.Range("H139").Value = 255728.75
$ws.Range("J139").Value = 255728.75
$ws.Range("L139").Value = 255728.75
$ws.Range("N139").Value = -266008.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row update (source hunk #10)
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
# Row update (source hunk #11)
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row update (source hunk #12)
$ws.Range("H16").Value = 2747.5
$ws.Range("I16").Value = 1992.3334
$ws.Range("K16").Value = 1992.3334
$ws.Range("M16").Value = -1705.3334
# Row update (source hunk #13)
$ws.Range("H31").Value = 7918.5
$ws.Range("I31").Value = 8405.823
$ws.Range("K31").Value = 8405.823
$ws.Range("M31").Value = -8110.823
# Row update (source hunk #14)
$ws.Range("H34").Value = 7918.5
$ws.Range("I34").Value = 8405.823
$ws.Range("K34").Value = 8405.823
$ws.Range("M34").Value = -8203.823
# Row update (source hunk #15)
$ws.Range("H58").Value = 3002.3572
$ws.Range("I58").Value = 3007.9
$ws.Range("J58").Value = 2988.5
$ws.Range("K58").Value = 3007.9
$ws.Range("L58").Value = 2988.5
$ws.Range("M58").Value = -2804.9
$ws.Range("N58").Value = -3394.5
# Row update (source hunk #16)
$ws.Range("H68").Value = 54990
$ws.Range("J68").Value = 54990
$ws.Range("L68").Value = 54990
$ws.Range("N68").Value = -56488
# Row update (source hunk #17)
$ws.Range("H71").Value = 54990
$ws.Range("J71").Value = 54990
$ws.Range("L71").Value = 164970
$ws.Range("N71").Value = -172458
# Row update (source hunk #18)
$ws.Range("H113").Value = 2747.5
$ws.Range("I113").Value = 1992.3334
$ws.Range("K113").Value = 1992.3334
$ws.Range("M113").Value = 177.6666
# Row update (source hunk #19)
$ws.Range("H136").Value = 3002.3572
$ws.Range("I136").Value = 3007.9
$ws.Range("J136").Value = 2988.5
$ws.Range("K136").Value = 9023.700000000001
$ws.Range("L136").Value = 8965.5
$ws.Range("M136").Value = -6473.700000000001
$ws.Range("N136").Value = -14065.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row update (source hunk #20)
$ws.Range("H2").Value = 70.42856999999999
$ws.Range("I2").Value = 63.444443
$ws.Range("K2").Value = 380.666658
$ws.Range("M2").Value = -267.666658
# Row update (source hunk #21)
$ws.Range("H19").Value = 300
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 900
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1248
# Row update (source hunk #22)
$ws.Range("H131").Value = 10872306
$ws.Range("J131").Value = 1964.3182
$ws.Range("L131").Value = 5892.9546
$ws.Range("N131").Value = -15972.9546

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row update (source hunk #23)
$ws.Range("H122").Value = 14058.5
$ws.Range("I122").Value = 8631
$ws.Range("J122").Value = 25999
$ws.Range("K122").Value = 25893
$ws.Range("L122").Value = 77997
$ws.Range("M122").Value = -23443
$ws.Range("N122").Value = -82897
# Row update (source hunk #24)
$ws.Range("H132").Value = 4628.575
$ws.Range("I132").Value = 4901.394
$ws.Range("J132").Value = 3342.4285
$ws.Range("K132").Value = 14704.182
$ws.Range("L132").Value = 10027.2855
$ws.Range("M132").Value = -12174.182
$ws.Range("N132").Value = -15087.2855

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row update (source hunk #25)
$ws.Range("H22").Value = 14306.667
$ws.Range("J22").Value = 1733.3334
$ws.Range("L22").Value = 1733.3334
$ws.Range("N22").Value = -2323.3334
# Row update (source hunk #26)
$ws.Range("H27").Value = 14306.667
$ws.Range("J27").Value = 1733.3334
$ws.Range("L27").Value = 1733.3334
$ws.Range("N27").Value = -1947.3334
# Row update (source hunk #27)
$ws.Range("H136").Value = 4674.4287
$ws.Range("I136").Value = 3950.7646
$ws.Range("J136").Value = 7750
$ws.Range("K136").Value = 11852.2938
$ws.Range("L136").Value = 23250
$ws.Range("M136").Value = -9302.293799999999
$ws.Range("N136").Value = -28350
# Row update (source hunk #28)
$ws.Range("H141").Value = 97226.5
$ws.Range("J141").Value = 97226.5
$ws.Range("L141").Value = 97226.5
$ws.Range("N141").Value = -107586.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row update (source hunk #29)
$ws.Range("H5").Value = 16676667
$ws.Range("I5").Value = 25000000
$ws.Range("J5").Value = 30000
$ws.Range("K5").Value = 25000000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = -24999888
$ws.Range("N5").Value = -30224
# Row update (source hunk #30)
$ws.Range("H132").Value = 8969.289000000001
$ws.Range("J132").Value = 5151.077
$ws.Range("L132").Value = 15453.231
$ws.Range("N132").Value = -20513.231
# Row update (source hunk #31)
$ws.Range("H136").Value = 489651.4
$ws.Range("I136").Value = 502220.8
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 1506662.4
$ws.Range("L136").Value = 300000
$ws.Range("M136").Value = -1504112.4
$ws.Range("N136").Value = -305100

Write-Host "All updates applied."